# Comment from python script
# Fill in column D ("T" score column) on the "Ninja" sheet for every
# student row (2-29). Every student gets a 1 except row 28, which gets 0.
$wb = $excel.ActiveWorkbook
$ninja = $wb.Worksheets.Item("Ninja")

$dValues = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1;
    19 = 1; 20 = 1; 21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1; 26 = 1;
    27 = 1; 28 = 0; 29 = 1
}

foreach ($row in 2..29) {
    $ninja.Cells.Item($row, 4).Value = $dValues[$row]
}

# The Senador sheet was the active/selected sheet before this edit;
# move off of it and leave its cursor parked at C2.
$senador = $wb.Worksheets.Item("Senador")
$senador.Activate() | Out-Null
$senador.Range("C2").Select() | Out-Null

# Finish on the Ninja sheet, which becomes the active tab, with the
# cursor left on the last-edited cell.
$ninja.Activate() | Out-Null
$ninja.Range("D26").Select() | Out-Null
